# implemented mixed languages on activity executions
# Adds a new "Gemischte Durchfuehrung" (mixed execution) column (H) to the
# activity-execution upload template, mirroring the existing "Mit Transport"
# (G) column's header/values/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column H, row 1 - same row as the other headers.
$ws.Cells.Item(1, 8).Value = "Gemischte Durchfuehrung"

# Sample data rows, matching the "ja" values used in column G (Mit Transport).
$ws.Cells.Item(2, 8).Value = "ja"
$ws.Cells.Item(3, 8).Value = "ja"

# Row 3's data cells carry the bold-ish "s=3" style; copy that formatting
# from G3 (same row, prior column) onto the new H3 cell so it matches.
$ws.Cells.Item(3, 7).Copy()
$ws.Cells.Item(3, 8).PasteSpecial(-4122)

# Mirror the author's final cursor position after typing the new header.
$ws.Range("H4").Select() | Out-Null
